$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '27.317.45'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  -3.71%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.850.00'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  -5.23%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.002'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.36%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '322.17'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.86%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '1.002'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -0.30%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4493'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -5.70%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3829'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -4.76%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '48.34'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -10.07%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.07848'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -6.88%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '1.013'
$ws.Range('D11').Style = 'Normal'

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '21.26'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -4.55%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '1.837.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -6.04%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.851'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -4.67%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '7.122'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.50%  '

# Row 16
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '1.002'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -0.56%  '

# Row 17
$ws.Range('B17').Value = 'ShibaInu'
$ws.Range('C17').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.00001031'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -3.47%  '

# Row 18
$ws.Range('B18').Value = 'Litecoin'
$ws.Range('C18').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '85.55'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -5.55%  '

# Row 19
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '0.06494'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -1.34%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.00'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -7.55%  '

# Row 21
$ws.Range('E21').Value = '  -0.37%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.464'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -5.98%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '27.317.92'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -3.77%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '10.79'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -5.72%  '

# Row 25
$ws.Range('E25').Value = '  -0.75%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.065.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -5.54%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '151.22'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -2.09%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '19.35'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -4.09%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '5.528'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -6.44%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '2.046'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -4.82%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '119.79'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.03%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.09311'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -3.10%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '1.471'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +1.62%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.9322'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -4.48%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '3.597'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -1.66%  '

# Row 36
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.253'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -6.01%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.02218'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -4.53%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05966'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -3.85%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.201'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -3.74%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '8.252'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -7.86%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.001'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.39%  '

# Row 42
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.5895'
$ws.Range('D42').Style = 'Normal'

# Row 43
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.1848'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -3.40%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.18'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -8.31%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.249'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -7.11%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.5642'
$ws.Range('D46').Style = 'Normal'

# Row 47
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '12.11'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -6.41%  '

# Row 48
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.926'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -6.18%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '3.359'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -0.72%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06875'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.10%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '108.04'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.27%  '

